$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4872029006772032
$ws.Range("C2").Value = 0.9902987556130791

$ws.Range("B3").Value = 0.2379542753380778
$ws.Range("C3").Value = 0.9953491777044162

$ws.Range("B4").Value = 0.3461213788481244
$ws.Range("C4").Value = 0.9933413943208247

$ws.Range("B5").Value = 0.4131214985353385
$ws.Range("C5").Value = 0.9918538096626445
